$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs",          "Ngf", "Ngfr", "ECs",              3, 1,               1.808798,            5.426394,  0.36666838522954,    0.4612006375427627, 2, 0.6666666666666666, 0.08108666666666665, 0.24326,   0.02056149724823249, 0.02939174694341418, 0.1466694004933333,  1.32002460444,       0.007539250993911038, 0.01355549242879817),
    @("FAPs",          "Ngf", "Ngfr", "FAPs",             3, 1,               1.808798,            5.426394,  0.36666838522954,    0.4612006375427627, 2, 0.6666666666666666, 0.090517,             0.271551,  0.02295278771378271, 0.03280999043916413, 0.163726968566,      1.473542717094,      0.00841606160752913,  0.01513198850831445),
    @("FAPs",          "Ngf", "Ngfr", "Inflammatory-Mac", 3, 1,               1.808798,            5.426394,  0.36666838522954,    0.4612006375427627, 2, 0.6666666666666666, 0.2176356666666667,  0.652907,  0.05518681856388939, 0.07888710565478799, 0.3936589585953334,  3.542930627358,      0.02023526164877693,  0.0363827834218915),
    @("FAPs",          "Ngf", "Ngfr", "MuSCs",            3, 1,               1.808798,            5.426394,  0.36666838522954,    0.4612006375427627, 2, 1,                   3.5543775,            7.108755,  0.9012988964740954,  0.8589111569626339,  6.429150913245,      38.57490547947,      0.3304778109793229,   0.3961303731837587),
    @("MuSCs",         "Ngf", "Ngfr", "ECs",              2, 1,               3.033388,            6.066776,  0.6149097244328354,  0.5156280504196953, 2, 0.6666666666666666, 0.08108666666666665, 0.24326,   0.02056149724823249, 0.02939174694341418, 0.2459673216266666,  1.47580392976,       0.01264346460683714,  0.01515520917486169),
    @("MuSCs",         "Ngf", "Ngfr", "FAPs",             2, 1,               3.033388,            6.066776,  0.6149097244328354,  0.5156280504196953, 2, 0.6666666666666666, 0.090517,             0.271551,  0.02295278771378271, 0.03280999043916413, 0.274573181596,      1.647439089576,      0.01411389236804749,  0.01691775140443504),
    @("MuSCs",         "Ngf", "Ngfr", "Inflammatory-Mac", 2, 1,               3.033388,            6.066776,  0.6149097244328354,  0.5156280504196953, 2, 0.6666666666666666, 0.2176356666666667,  0.652907,  0.05518681856388939, 0.07888710565478799, 0.6601734196386667,  3.961040517832,      0.03393491139544611,  0.04067640449203085),
    @("MuSCs",         "Ngf", "Ngfr", "MuSCs",            2, 1,               3.033388,            6.066776,  0.6149097244328354,  0.5156280504196953, 2, 1,                   3.5543775,            7.108755,  0.9012988964740954,  0.8589111569626339,  10.78180605597,      43.12722422388,      0.5542174560625047,   0.4428786853483678),
    @("Resolving-Mac", "Ngf", "Ngfr", "ECs",              1, 0.3333333333333333, 0.09087633333333334, 0.272629, 0.01842189033762463, 0.02317131203754203, 2, 0.6666666666666666, 0.08108666666666665, 0.24326,   0.02056149724823249, 0.02939174694341418, 0.007368858948888888, 0.06631973054,      0.0003787816474843095, 0.0006810453397543222),
    @("Resolving-Mac", "Ngf", "Ngfr", "FAPs",             1, 0.3333333333333333, 0.09087633333333334, 0.272629, 0.01842189033762463, 0.02317131203754203, 2, 0.6666666666666666, 0.090517,             0.271551,  0.02295278771378271, 0.03280999043916413, 0.008225853064333334, 0.07403267757900001, 0.0004228337382060829, 0.0007602505264146427),
    @("Resolving-Mac", "Ngf", "Ngfr", "Inflammatory-Mac", 1, 0.3333333333333333, 0.09087633333333334, 0.272629, 0.01842189033762463, 0.02317131203754203, 2, 0.6666666666666666, 0.2176356666666667,  0.652907,  0.05518681856388939, 0.07888710565478799, 0.01977793138922222, 0.178001382503,     0.001016645519666357, 0.001827917740865639),
    @("Resolving-Mac", "Ngf", "Ngfr", "MuSCs",            1, 0.3333333333333333, 0.09087633333333334, 0.272629, 0.01842189033762463, 0.02317131203754203, 2, 1,                   3.5543775,            7.108755,  0.9012988964740954,  0.8589111569626339,  0.3230087944825,    1.938052766895,     0.01660362943226788,  0.01990209843050744)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($r, $col).Value = $row[$c]
    }
}
